$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 9478.134
$ws.Range("I86").Value = 11355.637
$ws.Range("J86").Value = 4315
$ws.Range("K86").Value = 11355.637
$ws.Range("L86").Value = 4315
$ws.Range("M86").Value = -10232.637
$ws.Range("N86").Value = -6561
$ws.Range("H89").Value = 9478.134
$ws.Range("I89").Value = 11355.637
$ws.Range("J89").Value = 4315
$ws.Range("K89").Value = 56778.185
$ws.Range("L89").Value = 21575
$ws.Range("M89").Value = -51162.185
$ws.Range("N89").Value = -32807
$ws.Range("H98").Value = 2041.4117
$ws.Range("I98").Value = 978.9286
$ws.Range("J98").Value = 6999.6665
$ws.Range("K98").Value = 978.9286
$ws.Range("L98").Value = 6999.6665
$ws.Range("M98").Value = 519.0714
$ws.Range("N98").Value = -9995.666499999999
$ws.Range("H122").Value = 2041.4117
$ws.Range("I122").Value = 978.9286
$ws.Range("J122").Value = 6999.6665
$ws.Range("K122").Value = 2936.7858
$ws.Range("L122").Value = 20998.9995
$ws.Range("M122").Value = -486.7857999999997
$ws.Range("N122").Value = -25898.9995
$ws.Range("M125").Value = -4740
$ws.Range("H125").Value = 1854004.6
$ws.Range("I125").Value = 800
$ws.Range("K125").Value = 7200
$ws.Range("H129").Value = 899.4
$ws.Range("I129").Value = 404.78946
$ws.Range("J129").Value = 1015.41974
$ws.Range("K129").Value = 1214.36838
$ws.Range("L129").Value = 3046.25922
$ws.Range("M129").Value = 3785.63162
$ws.Range("N129").Value = -13046.25922
$ws.Range("H132").Value = 2687.3257
$ws.Range("I132").Value = 1554.5278
$ws.Range("J132").Value = 8513.143
$ws.Range("K132").Value = 4663.5834
$ws.Range("L132").Value = 25539.429
$ws.Range("M132").Value = -2133.5834
$ws.Range("N132").Value = -30599.429
$ws.Range("H137").Value = 2728306.8
$ws.Range("I137").Value = 1087934.6
$ws.Range("J137").Value = 11112431
$ws.Range("K137").Value = 3263803.8
$ws.Range("L137").Value = 33337293
$ws.Range("M137").Value = -3261253.8
$ws.Range("N137").Value = -33342393

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2600
$ws.Range("I88").Value = 2500
$ws.Range("J88").Value = 2666.6667
$ws.Range("K88").Value = 2500
$ws.Range("L88").Value = 2666.6667
$ws.Range("M88").Value = -2094
$ws.Range("N88").Value = -3478.6667
$ws.Range("H91").Value = 2600
$ws.Range("I91").Value = 2500
$ws.Range("J91").Value = 2666.6667
$ws.Range("K91").Value = 2500
$ws.Range("L91").Value = 2666.6667
$ws.Range("M91").Value = -1096
$ws.Range("N91").Value = -5474.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1349.5652
$ws.Range("I99").Value = 777
$ws.Range("J99").Value = 5166.6665
$ws.Range("K99").Value = 777
$ws.Range("L99").Value = 5166.6665
$ws.Range("M99").Value = 721
$ws.Range("N99").Value = -8162.6665
$ws.Range("H107").Value = 8002.625
$ws.Range("I107").Value = 8836.833000000001
$ws.Range("J107").Value = 5500
$ws.Range("K107").Value = 8836.833000000001
$ws.Range("L107").Value = 5500
$ws.Range("M107").Value = -6916.833000000001
$ws.Range("N107").Value = -9340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 1975.6666
$ws.Range("I19").Value = 1975.6666
$ws.Range("K19").Value = 1975.6666
$ws.Range("M19").Value = -1805.6666
$ws.Range("H24").Value = 1975.6666
$ws.Range("I24").Value = 1975.6666
$ws.Range("K24").Value = 1975.6666
$ws.Range("M24").Value = -1805.6666
$ws.Range("H58").Value = 904.26666
$ws.Range("I58").Value = 907.8378
$ws.Range("J58").Value = 887.75
$ws.Range("K58").Value = 907.8378
$ws.Range("L58").Value = 887.75
$ws.Range("M58").Value = -704.8378
$ws.Range("N58").Value = -1293.75
$ws.Range("H62").Value = 2735.7144
$ws.Range("I62").Value = 2542.8572
$ws.Range("J62").Value = 2928.5715
$ws.Range("K62").Value = 2542.8572
$ws.Range("L62").Value = 2928.5715
$ws.Range("M62").Value = -1918.8572
$ws.Range("N62").Value = -4176.5715
$ws.Range("H65").Value = 2735.7144
$ws.Range("I65").Value = 2542.8572
$ws.Range("J65").Value = 2928.5715
$ws.Range("K65").Value = 12714.286
$ws.Range("L65").Value = 14642.8575
$ws.Range("M65").Value = -9594.286
$ws.Range("N65").Value = -20882.8575
$ws.Range("H122").Value = 1590.6666
$ws.Range("I122").Value = 1328.8572
$ws.Range("J122").Value = 2507
$ws.Range("K122").Value = 3986.5716
$ws.Range("L122").Value = 7521
$ws.Range("M122").Value = -1536.5716
$ws.Range("N122").Value = -12421
$ws.Range("H134").Value = 5475.2095
$ws.Range("I134").Value = 6095.4443
$ws.Range("K134").Value = 18286.3329
$ws.Range("M134").Value = -15751.3329
$ws.Range("H136").Value = 904.26666
$ws.Range("I136").Value = 907.8378
$ws.Range("J136").Value = 887.75
$ws.Range("K136").Value = 2723.5134
$ws.Range("L136").Value = 2663.25
$ws.Range("M136").Value = -173.5133999999998
$ws.Range("N136").Value = -7763.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 682.23914
$ws.Range("I113").Value = 586.0833
$ws.Range("J113").Value = 716.17645
$ws.Range("K113").Value = 1758.2499
$ws.Range("L113").Value = 2148.52935
$ws.Range("M113").Value = 411.7501
$ws.Range("N113").Value = -6488.529350000001
$ws.Range("M116").ClearContents() | Out-Null
$ws.Range("H116").Value = 7000
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 7000
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 21000
$ws.Range("N116").Value = -27884
$ws.Range("H125").Value = 2534.8286
$ws.Range("I125").Value = 1079.75
$ws.Range("J125").Value = 2722.5806
$ws.Range("K125").Value = 3239.25
$ws.Range("L125").Value = 8167.7418
$ws.Range("M125").Value = 1680.75
$ws.Range("N125").Value = -18007.7418
$ws.Range("H129").Value = 1298.138
$ws.Range("I129").Value = 733.0769
$ws.Range("J129").Value = 1757.25
$ws.Range("K129").Value = 2199.2307
$ws.Range("L129").Value = 5271.75
$ws.Range("M129").Value = 2800.7693
$ws.Range("N129").Value = -15271.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N7").ClearContents() | Out-Null
$ws.Range("H7").Value = 15000000
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N8").ClearContents() | Out-Null
$ws.Range("H8").Value = 15000000
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1655.2759
$ws.Range("I136").Value = 1406.12
$ws.Range("K136").Value = 4218.36
$ws.Range("M136").Value = -1668.36

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M11").ClearContents() | Out-Null
$ws.Range("H11").Value = 3797502.5
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 3797502.5
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 3797502.5
$ws.Range("N11").Value = -3797786.5
$ws.Range("H27").Value = 39800
$ws.Range("J27").Value = 39800
$ws.Range("L27").Value = 39800
$ws.Range("N27").Value = -39938
$ws.Range("N115").Value = -41134
$ws.Range("H115").Value = 38000
$ws.Range("J115").Value = 38000
$ws.Range("L115").Value = 38000
$ws.Range("H136").Value = 1577.7222
$ws.Range("I136").Value = 1666.4419
$ws.Range("J136").Value = 1230.909
$ws.Range("K136").Value = 4999.3257
$ws.Range("L136").Value = 3692.727
$ws.Range("M136").Value = -2449.3257
$ws.Range("N136").Value = -8792.727000000001
